# "chore: adapt column header formatting to respective input file names"
#
# The AHB-diff sheet compared two format-versions of a message type using
# generic "_old" / "_new" column-header suffixes. This edit renames those
# headers to spell out the actual format versions being compared
# (FV2310 -> FV2404), wraps the data range in a proper Excel Table so the
# generated file is easier to filter/inspect, and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old header text -> new header text, keyed by column letter.
$headerRenames = [ordered]@{
    "A1" = "Segmentname_FV2310"
    "B1" = "Segmentgruppe_FV2310"
    "C1" = "Segment_FV2310"
    "D1" = "Datenelement_FV2310"
    "E1" = "Segment ID_FV2310"
    "F1" = "Code_FV2310"
    "G1" = "Qualifier_FV2310"
    "H1" = "Beschreibung_FV2310"
    "I1" = "Bedingungsausdruck_FV2310"
    "J1" = "Bedingung_FV2310"
    "K1" = "diff"
    "L1" = "Segmentname_FV2404"
    "M1" = "Segmentgruppe_FV2404"
    "N1" = "Segment_FV2404"
    "O1" = "Datenelement_FV2404"
    "P1" = "Segment ID_FV2404"
    "Q1" = "Code_FV2404"
    "R1" = "Qualifier_FV2404"
    "S1" = "Beschreibung_FV2404"
    "T1" = "Bedingungsausdruck_FV2404"
    "U1" = "Bedingung_FV2404"
}

foreach ($cellRef in $headerRenames.Keys) {
    $ws.Range($cellRef).Value = $headerRenames[$cellRef]
}

# Wrap the whole used range in a real Excel Table (ListObject) so the
# headers become filterable columns, named "Table1" like in the source
# workbook. Header names are picked up from the cells we just renamed.
$dataRange = $ws.Range("A1:U62")
$listObj = $ws.ListObjects.Add(1, $dataRange, 0, 1)
$listObj.Name = "Table1"

# Freeze the header row (split below row 1) so it stays visible on scroll.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
